# 1) Update the letter date.
$d = $word.ActiveDocument
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2) Split the mailing address into two lines: street on its own paragraph,
#    city/state/zip on a new paragraph right after it.
$addrRange = $d.Content
$addrRange.Find.Execute("2564 Greenrock Road, Milpitas CA 95035", $true, $false, $false,
                         $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($addrRange.Find.Found) {
    $addrRange.Text = "2564 Greenrock Road"
    $addrRange.Collapse(0)
    $addrRange.InsertParagraphAfter()
    $addrRange.Collapse(0)
    $addrRange.Move(1, 1) | Out-Null
    $addrRange.InsertAfter("Milpitas, CA 95035")
}

# 3) Remove the now-redundant blank "NoSpacing" paragraph that used to sit
#    directly under "Board of Directors".
foreach ($p in @($d.Paragraphs)) {
    $t = $p.Range.Text
    if ($t -ne $null -and $t.Trim() -eq "" -and $p.Style.NameLocal -eq "No Spacing") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text.Contains("Board of Directors")) {
            $p.Range.Delete() | Out-Null
            break
        }
    }
}
